$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell - match style/format of existing header cells (G1 "sum")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column values (0/1 flags per row, rows 2-7)
$saveValues = @(0, 0, 0, 0, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
